$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Add the required "experimental" value ("true") next to the "Experimental" property row.
# Typing the bare word true/false always becomes a native Excel boolean, so instead we
# build it as a text formula and then paste-special as values, which leaves behind a
# plain text cell (type "s") containing the literal string "true" - matching the
# generator tool's output - instead of a boolean ("b") cell.
$cellB7 = $ws.Range("B7")
$cellB7.Formula = "=""true"""
$cellB7.Copy()
$cellB7.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# Update the "Date" property value
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"
